$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template rows for the two alternating fill styles already used in the sheet
# (style "1" = light fill, style "2" = darker fill). Copying number formats from an
# existing row keeps the original cellXfs/fills untouched (no new styles are added).
$style1Template = "A2:E2"
$style2Template = "A9:E9"

# Row 54
$ws.Range($style1Template).Copy()
$ws.Range("A54:E54").PasteSpecial(-4122)
$ws.Range("A54").Value = "label49"
$ws.Range("B54").Value = "caseID"
$ws.Range("C54").Value = "comboBox14"
$ws.Range("D54").Value = 13
$ws.Range("E54").Value = "trans/linelineType/insert"

# Row 55
$ws.Range($style1Template).Copy()
$ws.Range("A55:E55").PasteSpecial(-4122)
$ws.Range("A55").Value = "label50"
$ws.Range("B55").Value = "line"
$ws.Range("C55").Value = "comboBox15"
$ws.Range("D55").Value = 13
$ws.Range("E55").Value = "trans/linelineType/insert"

# Row 56
$ws.Range($style1Template).Copy()
$ws.Range("A56:E56").PasteSpecial(-4122)
$ws.Range("A56").Value = "label51"
$ws.Range("B56").Value = "lineType"
$ws.Range("C56").Value = "comboBox16"
$ws.Range("D56").Value = 13
$ws.Range("E56").Value = "trans/linelineType/insert"

# Row 57
$ws.Range($style1Template).Copy()
$ws.Range("A57:E57").PasteSpecial(-4122)
$ws.Range("A57").Value = "button14"
$ws.Range("B57").Value = "Submit"
$ws.Range("C57").Value = "-"
$ws.Range("D57").Value = 13
$ws.Range("E57").Value = "trans/linelineType/insert"

# Row 58
$ws.Range($style1Template).Copy()
$ws.Range("A58:E58").PasteSpecial(-4122)
$ws.Range("A58").Value = "button15"
$ws.Range("B58").Value = "Clear"
$ws.Range("C58").Value = "-"
$ws.Range("D58").Value = 13
$ws.Range("E58").Value = "trans/linelineType/insert"

# Row 59
$ws.Range($style2Template).Copy()
$ws.Range("A59:E59").PasteSpecial(-4122)
$ws.Range("A59").Value = "label52"
$ws.Range("B59").Value = "case"
$ws.Range("C59").Value = "comboBox17"
$ws.Range("D59").Value = 14
$ws.Range("E59").Value = "trans/lossozine/insert"

# Row 60
$ws.Range($style2Template).Copy()
$ws.Range("A60:E60").PasteSpecial(-4122)
$ws.Range("A60").Value = "label53"
$ws.Range("B60").Value = "lossZOne"
$ws.Range("C60").Value = "textBox14"
$ws.Range("D60").Value = 14
$ws.Range("E60").Value = "trans/lossozine/insert"

# Row 61
$ws.Range($style2Template).Copy()
$ws.Range("A61:E61").PasteSpecial(-4122)
$ws.Range("A61").Value = "label54"
$ws.Range("B61").Value = "sequencialNumber"
$ws.Range("C61").Value = "textBox24"
$ws.Range("D61").Value = 14
$ws.Range("E61").Value = "trans/lossozine/insert"

# Row 62
$ws.Range($style2Template).Copy()
$ws.Range("A62:E62").PasteSpecial(-4122)
$ws.Range("A62").Value = "label55"
$ws.Range("B62").Value = "Description"
$ws.Range("C62").Value = "richTextBox4"
$ws.Range("D62").Value = 14
$ws.Range("E62").Value = "trans/lossozine/insert"

# Row 63
$ws.Range($style2Template).Copy()
$ws.Range("A63:E63").PasteSpecial(-4122)
$ws.Range("A63").Value = "button16"
$ws.Range("B63").Value = "submit"
$ws.Range("C63").Value = "-"
$ws.Range("D63").Value = 14
$ws.Range("E63").Value = "trans/lossozine/insert"

# Row 64
$ws.Range($style2Template).Copy()
$ws.Range("A64:E64").PasteSpecial(-4122)
$ws.Range("A64").Value = "label17"
$ws.Range("B64").Value = "clear"
$ws.Range("C64").Value = "-"
$ws.Range("D64").Value = 14
$ws.Range("E64").Value = "trans/lossozine/insert"

# Row 65
$ws.Range($style1Template).Copy()
$ws.Range("A65:E65").PasteSpecial(-4122)
$ws.Range("A65").Value = "label57"
$ws.Range("B65").Value = "caseID"
$ws.Range("C65").Value = "comboBox18"
$ws.Range("D65").Value = 15
$ws.Range("E65").Value = "trans/losszoneBus/insert"

# Row 66
$ws.Range($style1Template).Copy()
$ws.Range("A66:E66").PasteSpecial(-4122)
$ws.Range("A66").Value = "label58"
$ws.Range("B66").Value = "busID"
$ws.Range("C66").Value = "comboBox19"
$ws.Range("D66").Value = 15
$ws.Range("E66").Value = "trans/losszoneBus/insert"

# Row 67
$ws.Range($style1Template).Copy()
$ws.Range("A67:E67").PasteSpecial(-4122)
$ws.Range("A67").Value = "label59"
$ws.Range("B67").Value = "lossZoneID"
$ws.Range("C67").Value = "comboBox20"
$ws.Range("D67").Value = 15
$ws.Range("E67").Value = "trans/losszoneBus/insert"

# Row 68
$ws.Range($style1Template).Copy()
$ws.Range("A68:E68").PasteSpecial(-4122)
$ws.Range("A68").Value = "button18"
$ws.Range("B68").Value = "Submit"
$ws.Range("C68").Value = "-"
$ws.Range("D68").Value = 15
$ws.Range("E68").Value = "trans/losszoneBus/insert"

# Row 69
$ws.Range($style1Template).Copy()
$ws.Range("A69:E69").PasteSpecial(-4122)
$ws.Range("A69").Value = "button19"
$ws.Range("B69").Value = "Clear"
$ws.Range("C69").Value = "-"
$ws.Range("D69").Value = 15
$ws.Range("E69").Value = "trans/losszoneBus/insert"

$excel.CutCopyMode = 0

# Restore the selection to match the final cursor position from the edit session
$ws.Range("E67").Select()
